$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / update timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 00:22"

# --- Estados Unidos (row 4) refreshed totals ---
$ws.Range("B4").Value = 955491
$ws.Range("C4").Value = 30259
$ws.Range("D4").Value = 116186
$ws.Range("E4").Value = 785184
$ws.Range("G4").Value = 1928
$ws.Range("H4").Value = 54121

# --- Alemania (row 8) refreshed totals ---
$ws.Range("B8").Value = 156418
$ws.Range("C8").Value = 1419
$ws.Range("E8").Value = 40745
$ws.Range("G8").Value = 113
$ws.Range("H8").Value = 5873

# --- Belgica/Canada swap rank (Canada overtakes Belgica) ---
# Row 15 becomes Canada with the refreshed Canada totals
$ws.Range("A15").Value = "Canada"
$ws.Range("B15").Value = 45354
$ws.Range("C15").Value = 1466
$ws.Range("D15").Value = 16425
$ws.Range("E15").Value = 26464
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 163
$ws.Range("H15").Value = 2465

# Row 16 becomes Belgica, carrying the prior-day Belgica totals
$ws.Range("A16").Value = "Belgica"
$ws.Range("B16").Value = 45325
$ws.Range("C16").Value = 1032
$ws.Range("D16").Value = 10417
$ws.Range("E16").Value = 27991
$ws.Range("F16").Value = 934
$ws.Range("G16").Value = 238
$ws.Range("H16").Value = 6917

# --- Banglades/Colombia swap rank (Colombia overtakes Banglades) ---
# Row 50 becomes Colombia with the refreshed Colombia totals
$ws.Range("A50").Value = "Colombia"
$ws.Range("B50").Value = 5142
$ws.Range("C50").Value = 261
$ws.Range("D50").Value = 1067
$ws.Range("E50").Value = 3842
$ws.Range("F50").Value = 117
$ws.Range("G50").Value = 8
$ws.Range("H50").Value = 233

# Row 51 becomes Banglades, carrying the prior-day Banglades totals
$ws.Range("A51").Value = "Banglades"
$ws.Range("B51").Value = 4998
$ws.Range("C51").Value = 309
$ws.Range("D51").Value = 113
$ws.Range("E51").Value = 4745
$ws.Range("F51").Value = 1
$ws.Range("G51").Value = 9
$ws.Range("H51").Value = 140

# --- row 101 refreshed totals ---
$ws.Range("D101").Value = 442
$ws.Range("E101").Value = 146
